{"js": "// MHD2-259: Report template and related changes for reporting on 136 genes\n//\n// The document's single banner table had its header shading recolored:\n//   - table-level shading (w:tblPr/w:shd):  FFF2CC (theme accent4 tint) -> ECEAF2 (explicit)\n//   - title cell shading   (w:tcPr/w:shd):  E8E7EC                      -> ECEAF2\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document body.\");\n}\n\nconst table = tables.items[0];\n\n// Recolor the table's own shading (affects <w:tblPr><w:shd .../>).\ntable.shadingColor = \"#ECEAF2\";\n\n// Recolor the first (title) cell's shading (affects that cell's <w:tcPr><w:shd .../>).\nconst titleCell = table.getCell(0, 0);\ntitleCell.shadingColor = \"#ECEAF2\";\n\nawait context.sync();\n", "ps1": "# MHD2-259: Report template and related changes for reporting on 136 genes\n#\n# The document's single banner table had its header shading recolored:\n#   - table-level shading (w:tblPr/w:shd):  FFF2CC (theme accent4 tint) -> ECEAF2 (explicit)\n#   - title cell shading   (w:tcPr/w:shd):  E8E7EC                      -> ECEAF2\n\n$d = $word.ActiveDocument\n\n$table = $d.Tables.Item(1)\n\n# Recolor the table's own shading (affects <w:tblPr><w:shd .../>).\n$table.ShadingColor = \"#ECEAF2\"\n\n# Recolor the first (title) cell's shading (affects that cell's <w:tcPr><w:shd .../>).\n$cell = $table.Cell(1, 1)\n$cell.ShadingColor = \"#ECEAF2\"\n"}
